# Updated architecture image to correctly reflect Health Cloud's
# Self-Service Appointment Management.
#
# Helper: find a shape on a slide by its stable PowerPoint shape Id
# (shape Names are not unique in this deck, so Id-based lookup is used
# instead of Shapes.Item("Name")).
function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------
# 1. Widen the "Salesforce instance" rectangle (id 144) and slide its
#    icon (id 148) over so the box now spans the space previously
#    shared with the (now removed) duplicate "Healthcare provider"
#    box.
# ---------------------------------------------------------------
$rectSalesforce = Get-ShapeById $s.Shapes 144
$rectSalesforce.Left  = 149.826062992126
$rectSalesforce.Width = 283.46962029921264

$iconSalesforce = Get-ShapeById $s.Shapes 148
$iconSalesforce.Left = 149.64102962204726
$iconSalesforce.Top  = 29.165276590551183

# ---------------------------------------------------------------
# 2. Remove the duplicate "Healthcare provider" box (id 203) and its
#    icon (id 204) - the left-hand "Healthcare provider" swimlane is
#    being dropped from the diagram.
# ---------------------------------------------------------------
$dupRect = Get-ShapeById $s.Shapes 203
$dupRect.Delete()

$dupIcon = Get-ShapeById $s.Shapes 204
$dupIcon.Delete()

# ---------------------------------------------------------------
# 3. Reposition the "Patient" lane icon (id 205).
# ---------------------------------------------------------------
$patientIcon = Get-ShapeById $s.Shapes 205
$patientIcon.Left = 202.44299212598426
$patientIcon.Top  = 104.42338982677165

# ---------------------------------------------------------------
# 4. Update the label below that icon (id 206): it now reads
#    "Self-Service" / "Appointment" / "Management" instead of the
#    single line "Virtual appointment application", and moves/grows
#    to sit under the relocated icon.
# ---------------------------------------------------------------
$label = Get-ShapeById $s.Shapes 206
$label.TextFrame.TextRange.Text = "Self-Service`rAppointment`rManagement"
$label.Left = 147.71
$label.Top  = 148.44701387401577

# ---------------------------------------------------------------
# 5. Adjust the connector that links the icon (id 205) to the
#    Salesforce box (id 151) - it is now flipped vertically and
#    spans a shorter, repositioned run.
# ---------------------------------------------------------------
$connector1 = Get-ShapeById $s.Shapes 207
$connector1.VerticalFlip = -1
$connector1.Left   = 239.44299212598426
$connector1.Top    = 122.5732313464567
$connector1.Width  = 74.0855905511811
$connector1.Height = 0.3501584803149606

# ---------------------------------------------------------------
# 6. Lengthen the connector from the "Patient" icon (id 208) to the
#    repositioned icon (id 205).
# ---------------------------------------------------------------
$connector2 = Get-ShapeById $s.Shapes 210
$connector2.Width  = 144.7703937007874
$connector2.Height = 0.5880324960629921
